$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$colB = New-Object 'object[,]' 24,1
$colB[0,0] = 0.2007517892934061
$colB[1,0] = 0.1773979687425253
$colB[2,0] = 0.1630576458209703
$colB[3,0] = 0.1572140724234572
$colB[4,0] = 0.1562437789893778
$colB[5,0] = 0.1629788358104065
$colB[6,0] = 0.1926998711949182
$colB[7,0] = 0.2509570734005138
$colB[8,0] = 0.2937231989873226
$colB[9,0] = 0.313166991513782
$colB[10,0] = 0.3205278845686621
$colB[11,0] = 0.3189426837244866
$colB[12,0] = 0.3137726205771685
$colB[13,0] = 0.310605525760792
$colB[14,0] = 0.2924522421780864
$colB[15,0] = 0.2813126814608324
$colB[16,0] = 0.274904520717655
$colB[17,0] = 0.2727346747107902
$colB[18,0] = 0.2824986118147876
$colB[19,0] = 0.3152912530867411
$colB[20,0] = 0.3367109617565518
$colB[21,0] = 0.3252801412206452
$colB[22,0] = 0.2819624647449643
$colB[23,0] = 0.2352017568592544
$ws.Range("B2:B25").Value = $colB

$colC = New-Object 'object[,]' 24,1
$colC[0,0] = 0.04681665022438608
$colC[1,0] = 0.04563686530794797
$colC[2,0] = 0.04490388609007567
$colC[3,0] = 0.04460308260204471
$colC[4,0] = 0.04455300855323685
$colC[5,0] = 0.04489983781989793
$colC[6,0] = 0.04641167250347422
$colC[7,0] = 0.04930614197981953
$colC[8,0] = 0.05138723970470238
$colC[9,0] = 0.05232356054605702
$colC[10,0] = 0.05267657841026363
$colC[11,0] = 0.05260061933270066
$colC[12,0] = 0.05235263484077279
$colC[13,0] = 0.05220053418321413
$colC[14,0] = 0.05132583544740044
$colC[15,0] = 0.0507865390062392
$colC[16,0] = 0.05047537605648955
$colC[17,0] = 0.05036985599905819
$colC[18,0] = 0.05084404921775132
$colC[19,0] = 0.05242551629638825
$colC[20,0] = 0.0534500573044312
$colC[21,0] = 0.05290408529724999
$colC[22,0] = 0.05081805231728964
$colC[23,0] = 0.04853090527000603
$ws.Range("C2:C25").Value = $colC

$colD = New-Object 'object[,]' 24,1
$colD[0,0] = 0.03026650136908415
$colD[1,0] = 0.02852255846774199
$colD[2,0] = 0.02743928019702224
$colD[3,0] = 0.02699473200942748
$colD[4,0] = 0.02692072877599117
$colD[5,0] = 0.02743329738133582
$colD[6,0] = 0.02966780408092973
$colD[7,0] = 0.03394910246603189
$colD[8,0] = 0.03703157141753621
$colD[9,0] = 0.03841984229455164
$colD[10,0] = 0.03894350204611641
$colD[11,0] = 0.03883081429879809
$colD[12,0] = 0.03846296538353045
$colD[13,0] = 0.03823737933788607
$colD[14,0] = 0.03694056024778547
$colD[15,0] = 0.03614139920014026
$colD[16,0] = 0.03568043139031829
$colD[17,0] = 0.03552413160436174
$colD[18,0] = 0.03622660729035232
$colD[19,0] = 0.03857106744277417
$colD[20,0] = 0.04009134965431826
$colD[21,0] = 0.03928105409489291
$colD[22,0] = 0.03618808944258944
$colD[23,0] = 0.03280183929796721
$ws.Range("D2:D25").Value = $colD

$colF = New-Object 'object[,]' 24,1
$colF[0,0] = 0.7715222594088544
$colF[1,0] = 0.7677976138000417
$colF[2,0] = 0.7659605768555195
$colF[3,0] = 0.7653249729295908
$colF[4,0] = 0.7652262527364542
$colF[5,0] = 0.7659515475085072
$colF[6,0] = 0.7701445409289747
$colF[7,0] = 0.7819442100709537
$colF[8,0] = 0.7928066154587015
$colF[9,0] = 0.7982272655330007
$colF[10,0] = 0.80034903446969
$colF[11,0] = 0.7998889982797692
$colF[12,0] = 0.7984004391939408
$colF[13,0] = 0.7974976557930802
$colF[14,0] = 0.7924620218295786
$colF[15,0] = 0.7894956978067427
$colF[16,0] = 0.7878346495289748
$colF[17,0] = 0.7872799885536068
$colF[18,0] = 0.789806798896791
$colF[19,0] = 0.7988357887702762
$colF[20,0] = 0.8051395217436124
$colF[21,0] = 0.801738192663862
$colF[22,0] = 0.7896660120254921
$colF[23,0] = 0.7783678407808026
$ws.Range("F2:F25").Value = $colF

$colG = New-Object 'object[,]' 24,1
$colG[0,0] = 0.002432968549593463
$colG[1,0] = 0.002435432073174725
$colG[2,0] = 0.002437023900148529
$colG[3,0] = 0.002437692567007301
$colG[4,0] = 0.002437804807559507
$colG[5,0] = 0.002437032837048254
$colG[6,0] = 0.002433801571348831
$colG[7,0] = 0.002428090583965134
$colG[8,0] = 0.00242427182472807
$colG[9,0] = 0.002422615552596724
$colG[10,0] = 0.002421999929765803
$colG[11,0] = 0.002422132001384973
$colG[12,0] = 0.002422564673356263
$colG[13,0] = 0.002422831202396217
$colG[14,0] = 0.002424381688601152
$colG[15,0] = 0.002425353538329299
$colG[16,0] = 0.002425920138814075
$colG[17,0] = 0.002426113290356894
$colG[18,0] = 0.002425249295008724
$colG[19,0] = 0.002422437273653183
$colG[20,0] = 0.002420666874381601
$colG[21,0] = 0.002421605621212474
$colG[22,0] = 0.002425296398888708
$colG[23,0] = 0.002429569029591101
$ws.Range("G2:G25").Value = $colG

$colK = New-Object 'object[,]' 24,1
$colK[0,0] = 0.1710339568456334
$colK[1,0] = 0.1491899232329814
$colK[2,0] = 0.1357270085283773
$colK[3,0] = 0.1302285189044454
$colK[4,0] = 0.1293147735785425
$colK[5,0] = 0.1356529029091149
$colK[6,0] = 0.1635129140663309
$colK[7,0] = 0.2177264288051219
$colK[8,0] = 0.2572800361339205
$colK[9,0] = 0.2752096001141524
$colK[10,0] = 0.2819894808135359
$colK[11,0] = 0.2805297498345851
$colK[12,0] = 0.275767581445109
$colK[13,0] = 0.2728493430440153
$colK[14,0] = 0.256106973732841
$colK[15,0] = 0.2458194238918878
$colK[16,0] = 0.2398963444489226
$colK[17,0] = 0.2378898848169229
$colK[18,0] = 0.2469151709665454
$colK[19,0] = 0.277166611896007
$colK[20,0] = 0.2968811564415148
$colK[21,0] = 0.2863644743963505
$colK[22,0] = 0.2464198105557358
$colK[23,0] = 0.2031073748764101
$ws.Range("K2:K25").Value = $colK

$colM = New-Object 'object[,]' 24,1
$colM[0,0] = 0.8830009790449651
$colM[1,0] = 0.7864211113928548
$colM[2,0] = 0.7276523867524674
$colM[3,0] = 0.7038312175626658
$colM[4,0] = 0.6998832577335889
$colM[5,0] = 0.7273306176414991
$colM[6,0] = 0.849586131242674
$colM[7,0] = 1.093840125834618
$colM[8,0] = 1.276471696067603
$colM[9,0] = 1.36034162384216
$colM[10,0] = 1.392222276758872
$colM[11,0] = 1.385350710354274
$colM[12,0] = 1.362961997961719
$colM[13,0] = 1.349264235165819
$colM[14,0] = 1.27100711754774
$colM[15,0] = 1.223206392013182
$colM[16,0] = 1.195786578714745
$colM[17,0] = 1.186515172724967
$colM[18,0] = 1.228287159903431
$colM[19,0] = 1.369534764003674
$colM[20,0] = 1.462556663086772
$colM[21,0] = 1.412841828489334
$colM[22,0] = 1.225989953857507
$colM[23,0] = 1.027236340030839
$ws.Range("M2:M25").Value = $colM

$colO = New-Object 'object[,]' 24,1
$colO[0,0] = 2.668328567852228
$colO[1,0] = 2.669398983641855
$colO[2,0] = 2.671518554298189
$colO[3,0] = 2.672749496696582
$colO[4,0] = 2.672976056740566
$colO[5,0] = 2.671533669139137
$colO[6,0] = 2.668393823117356
$colO[7,0] = 2.673866028040379
$colO[8,0] = 2.685017499365387
$colO[9,0] = 2.691648342438157
$colO[10,0] = 2.694383987834897
$colO[11,0] = 2.693784814080715
$colO[12,0] = 2.691868898768433
$colO[13,0] = 2.690724625939538
$colO[14,0] = 2.68461556446249
$colO[15,0] = 2.681267352691378
$colO[16,0] = 2.679488155223169
$colO[17,0] = 2.678910909323008
$colO[18,0] = 2.681608597742354
$colO[19,0] = 2.692425546308243
$colO[20,0] = 2.700805014397531
$colO[21,0] = 2.696212655196518
$colO[22,0] = 2.681453866975033
$colO[23,0] = 2.671136376070137
$ws.Range("O2:O25").Value = $colO
